$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @{
    2 = @(3; 1; 1.384145666666667; 4.152437; 0.1014617184198512; 0.1334061399754118; 3; 1; 6.743825; 20.231475; 0.06175298453338176; 0.06809253940846915; 9.334436150508333; 84.009925354575; 0.006265563928311404; 0.00908396284360748)
    3 = @(3; 1; 1.384145666666667; 4.152437; 0.1014617184198512; 0.1334061399754118; 3; 1; 69.657865; 208.973595; 0.6378547871531949; 0.7033368922862506; 96.41663198900167; 867.7496879010149; 0.06471784280689157; 0.09382945990221067)
    4 = @(3; 1; 1.384145666666667; 4.152437; 0.1014617184198512; 0.1334061399754118; 3; 1; 1.263827666666667; 3.791483; 0.01157282852869501; 0.01276089388411081; 1.749321588230111; 15.743894294071; 0.001174199069499673; 0.001702381595715063)
    5 = @(3; 1; 1.384145666666667; 4.152437; 0.1014617184198512; 0.1334061399754118; 3; 1; 1.038909; 3.116727; 0.009513255668495419; 0.01048988549144043; 1.438001390411; 12.942012513699; 0.0009652312678929353; 0.001399415132197144)
    6 = @(3; 1; 1.384145666666667; 4.152437; 0.1014617184198512; 0.1334061399754118; 2; 1; 30.502036; 61.004072; 0.2793061441162328; 0.205319788929729; 42.21926095391067; 253.315565723464; 0.02833888134725558; 0.02739092050168143)
    7 = @(3; 1; 2.458038666666667; 7.374116000000001; 0.1801810554109116; 0.2369096391566985; 3; 1; 6.743825; 20.231475; 0.06175298453338176; 0.06809253940846915; 16.57658261123333; 149.1892435011; 0.01112671792799842; 0.0161317789405237)
    8 = @(3; 1; 2.458038666666667; 7.374116000000001; 0.1801810554109116; 0.2369096391566985; 3; 1; 69.657865; 208.973595; 0.6378547871531949; 0.7033368922862506; 171.2217256074467; 1540.99553046702; 0.114929348748165; 0.1666272893571294)
    9 = @(3; 1; 2.458038666666667; 7.374116000000001; 0.1801810554109116; 0.2369096391566985; 3; 1; 1.263827666666667; 3.791483; 0.01157282852869501; 0.01276089388411081; 3.106537272669778; 27.958835454028; 0.002085204458389773; 0.003023178765401614)
    10 = @(3; 1; 2.458038666666667; 7.374116000000001; 0.1801810554109116; 0.2369096391566985; 3; 1; 1.038909; 3.116727; 0.009513255668495419; 0.01048988549144043; 2.553678493148; 22.983106438332; 0.001714108446743342; 0.00248515498657224)
    11 = @(3; 1; 2.458038666666667; 7.374116000000001; 0.1801810554109116; 0.2369096391566985; 2; 1; 30.502036; 61.004072; 0.2793061441162328; 0.205319788929729; 74.97518390005868; 449.851103400352; 0.05032567582961499; 0.04864223710707161)
    12 = @(2; 1; 9.799864; 19.599728; 0.7183572261692373; 0.6296842208678898; 3; 1; 6.743825; 20.231475; 0.06175298453338176; 0.06809253940846915; 66.08856783979999; 396.5314070388; 0.04436070267707193; 0.04287679762433798)
    13 = @(2; 1; 9.799864; 19.599728; 0.7183572261692373; 0.6296842208678898; 3; 1; 69.657865; 208.973595; 0.6378547871531949; 0.7033368922862506; 682.63760353036; 4095.825621182159; 0.4582075955981383; 0.4428801430269106)
    14 = @(2; 1; 9.799864; 19.599728; 0.7183572261692373; 0.6296842208678898; 3; 1; 1.263827666666667; 3.791483; 0.01157282852869501; 0.01276089388411081; 12.38533925277066; 74.312035516624; 0.008313425000805562; 0.008035333522994138)
    15 = @(2; 1; 9.799864; 19.599728; 0.7183572261692373; 0.6296842208678898; 3; 1; 1.038909; 3.116727; 0.009513255668495419; 0.01048988549144043; 10.181166908376; 61.087001450256; 0.006833915953859142; 0.00660531537267105)
    16 = @(2; 1; 9.799864; 19.599728; 0.7183572261692373; 0.6296842208678898; 2; 1; 30.502036; 61.004072; 0.2793061441162328; 0.205319788929729; 298.915804523104; 1195.663218092416; 0.2006415869393622; 0.129286631320976)
}

foreach ($r in $rowsData.Keys) {
    $vals = $rowsData[$r]
    $col = 5
    foreach ($v in $vals) {
        $ws.Cells.Item([int]$r, $col).Value = $v
        $col = $col + 1
    }
}
